# Remove the 4 data rows that belong to DOI "10.1038:s41556-022-01079-4"
# (fig1..fig4), which were miscategorized -- "Change to others mistakes".
# All following rows shift up by 4 as a result.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14:K17").EntireRow.Delete() | Out-Null

# Column A (Fig Index) was re-sized to fit its (now shorter) longest entry.
$ws.Columns("A").ColumnWidth = 32

# Reflect the rows that were selected (and deleted) as the active selection,
# mirroring the author's on-screen state at save time.
$ws.Range("A14:XFD17").Select() | Out-Null
